$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert new row for the "Estimate Volumes" (EV) shortcut above the old row 84
$ws.Rows.Item(84).Insert()
$ws.Range("B84").Value = "EV"
$ws.Range("A84").Value = "Estimate Volumes"
$ws.Rows.Item(84).RowHeight = 17

# Update print area to extend by one row
$ws.PageSetup.PrintArea = '$A$1:$C$129'

# Update active cell / selection to reflect the edit location
$ws.Range("A85").Select() | Out-Null

Write-Host "done"
